$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns involved in the row-content rotation: A, B, D, E, F, G, H, M, Q, R
$cols = @("A", "B", "D", "E", "F", "G", "H", "M", "Q", "R")

# Capture the original ("before") values for rows 3, 5, 6, 7 so the rotation
# can be applied without clobbering source data while writing.
$orig = @{}
foreach ($r in @(3, 5, 6, 7)) {
    $orig[$r] = @{}
    foreach ($c in $cols) {
        $orig[$r][$c] = $ws.Range("$c$r").Value()
    }
}

# The row content rotates: new(3) = old(6), new(6) = old(5), new(5) = old(7), new(7) = old(3)
$mapping = @{
    3 = 6
    6 = 5
    5 = 7
    7 = 3
}

foreach ($target in $mapping.Keys) {
    $source = $mapping[$target]
    foreach ($c in $cols) {
        $ws.Range("$c$target").Value = $orig[$source][$c]
    }
}
